# Apply new "Applied Voltage [V]" (column M) values in 5 blocks of 20 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$voltages = @(0.5, 2.5, 4.5, 6.5, 8.5)
$blockSize = 20
$startRow = 2

for ($block = 0; $block -lt $voltages.Length; $block++) {
    $value = $voltages[$block]
    $firstRow = $startRow + ($block * $blockSize)
    $lastRow = $firstRow + $blockSize - 1
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 13).Value = $value
    }
}

# Update the sheet view: scroll so D11 is the top-left visible cell, and
# select M2:M101 (matching the saved selection state in the workbook).
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 4
$ws.Range("M2:M101").Select() | Out-Null
